# modification task et codesystem
# ajout demande par ROR et traitement du warning caseSensitive sur les codeSystems
#
# - Bump the "Date" metadata value (row 8, column B on the "Metadata" sheet)
#   to the new generation timestamp.
# - Fill in the previously-empty "Case Sensitive" value (row 14, column B)
#   with the text "true" (stored as a normal shared string, NOT as an Excel
#   boolean - Excel auto-coerces a literal Value of "true"/"false" into a
#   boolean cell, so we round-trip the text through a throw-away formula
#   cell + PasteSpecial(values) to force a text/string result instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Update the Date row -------------------------------------------------
$ws.Cells.Item(8, 2).Value = "2023-04-12T13:10:15+00:00"

# --- Fill in Case Sensitive = "true" (as text, preserving style) ---------
$caseSensitiveCell = $ws.Cells.Item(14, 2)

# Staging cell, well outside the used range, holds a formula whose result
# is the text string "true" (xlFormulaText), so copying it over preserves
# the destination's existing number format / style and writes a plain
# string value instead of a boolean.
$stage = $ws.Cells.Item(1, 5)
$stage.Formula = '="true"'
$stage.Copy()
$caseSensitiveCell.PasteSpecial(-4163)  # xlPasteValues
$stage.ClearContents()
